$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Edited" column (F) with the new translation drafts for each line.
$ws.Range("F2").Value = "-Do you want to go play together? I'm pretty free."
$ws.Range("F4").Value = "And, I'm given a scowl."
$ws.Range("F5").Value = "You always do things like that, don't you?!"
$ws.Range("F7").Value = "C-come on, that's not it, wait a second..."
$ws.Range("F8").Value = "Don't keep saying 'pervert' in a place this like, please...."
$ws.Range("F11").Value = "Somehow, the way she's acting so recklessly, it kind of reminds me of Yuki...."
$ws.Range("F12").Value = "Lowest of the low! You're the lowest! The worst! Just go somewhere else!"
$ws.Range("F13").Value = "C-come on, that's not it...! Ah, geez, forget it, I'm going home. Alright, see you later!"
$ws.Range("F14").Value = "I don't know what she'd tell me if I stayed here any longer. I quickly leave the place."

# Clear the lingering cell selection from F9 so the saved view matches a fresh A1 selection.
$ws.Range("A1").Select()
